# Fruta / hortaliza, semanal
# Permutes the per-row "weekly" fields (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg) among
# the data rows of the sheet, reshuffling which week's figures land on
# which row while leaving the fixed descriptive columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (values currently sitting in
# the source row get moved onto the destination row).
$map = @{
    2  = 15
    3  = 13
    4  = 11
    5  = 4
    6  = 14
    7  = 8
    8  = 9
    9  = 7
    11 = 5
    12 = 16
    13 = 2
    14 = 17
    15 = 3
    16 = 6
    17 = 12
}

$cols = @("D", "M", "N", "O", "P", "R", "S")

# Snapshot the original values for every affected column/row before
# writing anything back, since rows feed each other (a permutation, not
# a simple copy) and in-place writes would clobber values still needed
# as a source for another row.
$original = @{}
foreach ($row in $map.Values) {
    if (-not $original.ContainsKey($row)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$row").Value2
        }
        $original[$row] = $rowVals
    }
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
